$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update D2 with the corrected value
$ws.Range("D2").Value = 0.99660008614738449

# Rows 3-11: clear the "Image Name" (column A) values and zero out the
# LSB / LSB_pair / LSB-pair-ultar measurement columns (B:D)
for ($r = 3; $r -le 11; $r++) {
    $ws.Cells.Item($r, 1).ClearContents()
    $ws.Cells.Item($r, 2).Value = 0
    $ws.Cells.Item($r, 3).Value = 0
    $ws.Cells.Item($r, 4).Value = 0
}

$wb.Save()
